$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns, plus the swapped
# NEARProtocol/EnergySwap rows (44-45), to reflect the refreshed
# cryptos data pulled by the GitHub Actions job.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.944.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.904.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5025"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4043"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08252"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.904.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.351"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.181"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.53%  "
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06489"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.930"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.969.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.188"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.125.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.264"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1032"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.894"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.797"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.369"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06324"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2136"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.189"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6454"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.609"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.37%  "
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.187"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.75%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.132"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.68%  "
